$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.974.45"
$ws.Range("E2").Value = "  +0.20%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.638.75"
$ws.Range("E3").Value = "  -0.55%  "

$ws.Range("E4").Value = "  -0.73%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.03"
$ws.Range("E5").Value = "  -0.23%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5130"
$ws.Range("E6").Value = "  +0.84%  "

$ws.Range("E7").Value = "  -0.53%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2578"
$ws.Range("E8").Value = "  +0.03%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06353"
$ws.Range("E9").Value = "  -1.03%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.77"
$ws.Range("E10").Value = "  +0.32%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07755"
$ws.Range("E11").Value = "  -0.25%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.278"
$ws.Range("E12").Value = "  -0.83%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.642.43"
$ws.Range("E13").Value = "  -1.20%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5459"
$ws.Range("E14").Value = "  -0.19%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0₅7755"
$ws.Range("E15").Value = "  -1.80%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.36"
$ws.Range("E16").Value = "  -1.06%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "25.992.95"
$ws.Range("E17").Value = "  -0.06%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9998"
$ws.Range("E18").Value = "  -0.72%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "197.92"
$ws.Range("E19").Value = "  +0.31%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.446"
$ws.Range("E20").Value = "  +0.30%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.933"
$ws.Range("E21").Value = "  -1.05%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.082"
$ws.Range("E22").Value = "  +0.17%  "

$ws.Range("E23").Value = "  -0.66%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.931"
$ws.Range("E24").Value = "  +3.77%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "141.64"
$ws.Range("E25").Value = "  +0.53%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1233"
$ws.Range("E26").Value = "  +7.45%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.849"
$ws.Range("E27").Value = "  -0.66%  "

$ws.Range("E28").Value = "  -1.14%  "

$ws.Range("E29").Value = "  +0.07%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.04843"
$ws.Range("E30").Value = "  -3.45%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.283"
$ws.Range("E31").Value = "  +0.25%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.217"
$ws.Range("E32").Value = "  +0.43%  "

$ws.Range("E33").Value = "  -0.50%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.374"
$ws.Range("E34").Value = "  +0.16%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9145"
$ws.Range("E35").Value = "  +2.14%  "

$ws.Range("B36").Value = "MXToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.568"
$ws.Range("E36").Value = "  -0.95%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.5558"
$ws.Range("E37").Value = "  +0.09%  "

$ws.Range("B38").Value = "Maker"
$ws.Range("C38").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.133.71"
$ws.Range("E38").Value = "  +0.04%  "

$ws.Range("E39").Value = "  +0.44%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.001"
$ws.Range("E40").Value = "  -0.69%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.520"
$ws.Range("E41").Value = "  -1.88%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.584"
$ws.Range("E42").Value = "  -1.50%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8063"
$ws.Range("E43").Value = "  -1.29%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "99.42"
$ws.Range("E44").Value = "  -0.34%  "

$ws.Range("E45").Value = "  -3.40%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.780.94"
$ws.Range("E46").Value = "  -0.23%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4533"
$ws.Range("E47").Value = "  -0.19%  "

$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "55.10"
$ws.Range("E48").Value = "  -0.37%  "

$ws.Range("B49").Value = "Frax"
$ws.Range("C49").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.001"
$ws.Range("E49").Value = "  -0.66%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05216"
$ws.Range("E50").Value = "  +2.35%  "

$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.09548"
$ws.Range("E51").Value = "  +0.03%  "

